$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = "[name=`"Lens`"]'Is there anything else you would like to listen to?'`n"
$ws.Range("C34").Value = "[name=`"Lens`"]'Okay, Lens will look it up in the library.'`n"
$ws.Range("C39").Value = "[name=`"Lens`"]'Oh, hello, Kroos.'`n"
$ws.Range("C41").Value = "[name=`"Lens`"]'Ooh—'`n"
$ws.Range("C52").Value = "[name=`"Lens`"]'Ooh—'`n"
$ws.Range("C42").Value = "[name=`"Lens`"]'Miss Scene says she doesn't mind. Go ahead!'`n"
$ws.Range("C45").Value = "[name=`"Lens`"]'You look pretty tired.'`n"
$ws.Range("C48").Value = "[name=`"Lens`"]'Oh, I can hear her now. It's Fang.'`n"
$ws.Range("C54").Value = "[name=`"Lens`"]'Miss Scene says it sounds very interesting. Lens will therefore be providing all the assistance you need, Kroos.'`n"
$ws.Range("C56").Value = "[name=`"Lens`"]'However, I can't guarantee that Fang won't recognize you.'`n"
$ws.Range("C58").Value = "[name=`"Lens`"]'If you say so—'`n"
$ws.Range("C80").Value = "[name=`"Lens`"]'Congratulations, Kroos. It looks like you're still hidden.'`n"
$ws.Range("C83").Value = "[name=`"Lens`"]'Lens uses a revolutionary new camouflage technique that combines stealth and optical refraction technologies, among others, to satisfy Miss Scene's needs.'`n"
$ws.Range("C85").Value = "[name=`"Lens`"]'I'm sorry, Miss Scene. Lens forgot to put on music for you.'`n"
$ws.Range("C86").Value = "[name=`"Lens`"]'And without further ado—'`n"
$ws.Range("C90").Value = "[name=`"Lens`"]'Sure thing!'`n"
$ws.Range("C92").Value = "[name=`"Lens`"]'Miss Scene is brainstorming ideas for her next piece.'`n"
$ws.Range("C93").Value = "[name=`"Lens`"]'She may look like she's spacing out, but make no mistake. She's thinking.'`n"
$ws.Range("C94").Value = "[name=`"Lens`"]'And while she has her thinking cap on, it is Lens's job to play music and deliver food to her.'`n"
$ws.Range("C96").Value = "[name=`"Lens`"]'Let's see—Not according to Lens's records.'`n"
$ws.Range("C97").Value = "[name=`"Lens`"]'If anyone tried to, Lens would activate these camouflage techniques. No one would be able to find her.'`n"
$ws.Range("C100").Value = "[name=`"Lens`"]'I believe Miss Scene said something similar before.'`n"
$ws.Range("C101").Value = "[name=`"Lens`"]'Something along the lines of, 'I wish I could live the way Cautuses do—' Something like that.'`n"
$ws.Range("C108").Value = "[name=`"Lens`"]'What Miss Scene is saying is you can make your own decisions, and that is a joy in itself.'`n"
$ws.Range("C110").Value = "[name=`"Lens`"]'Let's use Lens as an analogy.'`n"
$ws.Range("C111").Value = "[name=`"Lens`"]'Lens is able to take photos very quickly. Through the use of my system's auxiliary operations, Lens is capable of snapping photos of running Kuranta and storing those photos in my album without issues.'`n"
$ws.Range("C112").Value = "[name=`"Lens`"]'But Miss Scene can't do that.'`n"
$ws.Range("C113").Value = "[name=`"Lens`"]'Even when she's photographing people, Miss Scene takes a very long time, almost like she's sketching.'`n"
$ws.Range("C114").Value = "[name=`"Lens`"]'That's why Miss Scene is a good scenery photographer. She's not good at it because she likes taking photos of landscapes.'`n"
$ws.Range("C115").Value = "[name=`"Lens`"]'She wants to take photos of everything, but, most of the time, she's only able to take photos of still objects. Only they have the patience to wait until she's done without moving.'`n"
$ws.Range("C116").Value = "[name=`"Lens`"]'And that's what Lens meant.'`n"
$ws.Range("C117").Value = "[name=`"Lens`"]'You don't like to move, but if you put your mind to it, you can still be pretty agile.'`n"
$ws.Range("C119").Value = "[name=`"Lens`"]'Both Lens and Click have edited operational records of your missions, Kroos.'`n"
$ws.Range("C120").Value = "[name=`"Lens`"]'You usually just hang out around the destination, but once an order comes, you always head towards the new destination quickly. That makes you no different from the other Cautuses.'`n"
$ws.Range("C123").Value = "[name=`"Lens`"]'Without Lens, Miss Scene wouldn't even be able to move to the next location.'`n"
$ws.Range("C124").Value = "[name=`"Lens`"]'Lens believes there's a difference between 'can't do' and 'don't want to do.''`n"
$ws.Range("C126").Value = "[name=`"Lens`"]''There are too many things that we can't do, so let's do the things that we can do the best we can.''`n"
$ws.Range("C127").Value = "[name=`"Lens`"]'This is something that Scene told Lens before.'`n"
$ws.Range("C128").Value = "[name=`"Lens`"]'Of course, Miss Scene may be able to do fewer things, but the things that she can do, she does with relative ease.'`n"
$ws.Range("C129").Value = "[name=`"Lens`"]'In your case, Lens supposes you probably focus on putting your energy where it matters because you know your limits?'`n"
$ws.Range("C131").Value = "[name=`"Lens`"]'Looks like Lens's deduction is correct.'`n"
$ws.Range("C133").Value = "[name=`"Lens`"]'Lens believes in you, Miss Kroos.'`n"
$ws.Range("C136").Value = "[name=`"Lens`"]'Oh, sorry about that. The reception here isn't very good.'`n"
$ws.Range("C137").Value = "[name=`"Lens`"]'The motivational r-routine that Lens s-selected to inspire you had a few e-errors.'`n"
$ws.Range("C139").Value = "[name=`"Lens`"]'It looks like y-you w-w-w—'`n"
$ws.Range("C188").Value = "Mayer's Room, a.k.a. 'Lutra Workshop'`n"
$ws.Range("C228").Value = "[name=`"Lens`"]'Why, hello, Kroos.'`n"
$ws.Range("C229").Value = "[name=`"Lens`"]'Oh, you're here too, Mayer. Welcome!'`n"
$ws.Range("C232").Value = "[name=`"Lens`"]'Lens actually rebooted not long after you left.'`n"
$ws.Range("C233").Value = "[name=`"Lens`"]'Lens sometimes runs into problems with certain functions. It isn't uncommon for Lens to reboot.'`n"
$ws.Range("C235").Value = "[name=`"Lens`"]'Oh, you mean what Miss Scene was trying to say? She was trying to tell you, 'It's okay.''`n"
$ws.Range("C236").Value = "[name=`"Lens`"]'She takes a liiiittle while to get her words out, so she might've given you the wrong impression.'`n"
$ws.Range("C240").Value = "[name=`"Lens`"]'Sorry about that, Miss Mayer. Lens is doing very well.'`n"
$ws.Range("C244").Value = "[name=`"Lens`"]'That's right.'`n"
$ws.Range("C249").Value = "[name=`"Lens`"]'Closure has already ordered replacement parts for Lens. It won't be long before Lens receives an upgrade.'`n"
$ws.Range("C253").Value = "[name=`"Mayer`"]She said it was some kind of 'advanced component used in extreme circumstances to be installed on a high-end mobile photography platform.'`n"
$ws.Range("C257").Value = "[name=`"Lens`"]'Thank you very much, Mayer.'`n"
$ws.Range("C258").Value = "[name=`"Lens`"]'But Lens is afraid it takes Miss Scene a while to get around too. Lens will walk with Miss Scene, so we likely won't be able to come with you straight to Closure's.'`n"
$ws.Range("C262").Value = "[name=`"Lens`"]'If that's what Miss Scene says…'`n"
$ws.Range("C264").Value = "[name=`"Lens`"]'Oh, Mayer, please don't push me around like that!'`n"
